$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: B12 becomes a true numeric value (was stored as text "3")
$ws.Range("B12").Value = 3

# Row 13 (new row) — copy of row 12's original annotation, but with an
# updated sentence_purpose / issue_type / id / source_file / text,
# and politeness_score kept as text "3" (matching the original diff).
$ws.Range("A13").Value = "Ruilin"
$ws.Range("B13").Value = "'3"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "无"
$ws.Range("D13").Value = "CRT"
$ws.Range("E13").Value = "RES"
$ws.Range("F13").Value = "77474e59-42ef-43e4-850b-a07d6b41a266"
$ws.Range("G13").Value = "Syg-YfWCW_annotated.xlsx"
$ws.Range("H13").Value = "You absolutely know this but you hide these results."
